$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write one data row (columns B..H) given a 1-based sheet row number
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- 1) "Haiti" moves from its old spot (row 151, right after Guadalupe) to
#        right after "Martinica" (row 146) -> becomes row 147, with updated
#        daily numbers. Birmania/Togo/Suazilandia/Guadalupe each shift down
#        one row (147->148, 148->149, 149->150, 150->151); Martinica (146)
#        and Gibraltar (152) stay put.
$ws.Cells.Item(151, 1).Value = "Guadalupe"
Set-Row 151 154 0 104 37 4 0 13

$ws.Cells.Item(150, 1).Value = "Suazilandia"
Set-Row 150 163 0 14 147 0 0 2

$ws.Cells.Item(149, 1).Value = "Togo"
Set-Row 149 173 20 89 73 0 1 11

$ws.Cells.Item(148, 1).Value = "Birmania"
Set-Row 148 180 2 72 102 0 0 6

$ws.Cells.Item(147, 1).Value = "Haiti"
Set-Row 147 182 31 17 150 0 3 15

# --- 2) "Islas Virgenes Britanicas" moves from right after "Butan" to right
#        after "Papua Nueva Guinea" -> rows 212/213 swap places.
$ws.Cells.Item(212, 1).Value = "Islas Virgenes Britanicas"
Set-Row 212 7 0 4 2 0 0 1

$ws.Cells.Item(213, 1).Value = "Butan"
Set-Row 213 7 0 5 2 0 0 0

# --- 3) Plain numeric updates (counts refreshed for the day) ---
# Estados Unidos
Set-Row 4 1367397 20088 256282 1030341 16514 737 80774

# Japon
Set-Row 36 15777 114 8127 7026 287 17 624

# Kazajistan
Set-Row 59 5090 115 1941 3118 31 0 31

# Venezuela
Set-Row 128 414 12 193 211 2 0 10

# Isla de Man (only "Casos criticos", column F, changes)
$ws.Cells.Item(131, 6).Value = 21

# Dominica (only columns D and E change)
$ws.Cells.Item(198, 4).Value = 15
$ws.Cells.Item(198, 5).Value = 1

Write-Host "Applied countries & provincias Spain update"
